# Weekly update: a new week's price record is inserted as row 21
# (pushing the existing rows 21-58 down to 22-59, and growing the
# used range from A1:R58 to A1:R59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 (shifts rows 21:58 down to 22:59).
$ws.Rows(21).Insert()

# Fill in the new row with this week's record.
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44571
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112031
$ws.Range("G21").Value = "Poroto verde"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 31000
$ws.Range("M21").Value = 30500
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1220
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
